# Logboek FotoSjaak - add the "week 50" log sheet (Donderdag 12-12-2013):
# "Heb de Login Class test script gemaakt" (commit: "Test script login gemaakt")

$wb = $excel.ActiveWorkbook

# The new week's log re-uses the exact same layout/styling as "week 48", so
# clone that sheet and drop it in right after "week 48" (i.e. right before
# "Totaal"), then rename it.
$ws48 = $wb.Worksheets.Item("week 48")
$ws48.Copy($null, $ws48)
$ws50 = $wb.Worksheets.Item(3)
$ws50.Name = "week 50"

# Row 7: the one logged activity for week 50 - Thursday 12 Dec 2013, 09:15-09:33.
$ws50.Range("A7").Value = "Donnderdag"
$ws50.Range("B7").Value = 41620
$ws50.Range("C7").Value = 0.38541666666666669
$ws50.Range("D7").Value = 0.3979166666666667
$ws50.Range("F7").Value = "Heb de Login Class test script gemaakt"

# The rest of the copied rows (8,9,11-15) carried over week 48's entries -
# clear those out since week 50 only has the single entry above.
$ws50.Range("C8:D9").ClearContents()
$ws50.Range("F8:F9").ClearContents()
$ws50.Range("A11:D11").ClearContents()
$ws50.Range("F11").ClearContents()
$ws50.Range("C12:D12").ClearContents()
$ws50.Range("F12").ClearContents()
$ws50.Range("C13:D13").ClearContents()
$ws50.Range("F13").ClearContents()
$ws50.Range("C14:D14").ClearContents()
$ws50.Range("F14").ClearContents()
$ws50.Range("C15:D15").ClearContents()
$ws50.Range("F15").ClearContents()

# Those rows were taller to wrap the (now removed) activity text - let them
# shrink back to their natural auto height.
$ws50.Range("A8:H9").EntireRow.AutoFit()
$ws50.Range("A11:H11").EntireRow.AutoFit()
$ws50.Range("A12:H12").EntireRow.AutoFit()
$ws50.Range("A13:H13").EntireRow.AutoFit()
$ws50.Range("A14:H14").EntireRow.AutoFit()
$ws50.Range("A15:H15").EntireRow.AutoFit()

# Make "week 50" the active sheet/selection, matching the saved view state.
$ws50.Activate()
$ws50.Range("F7").Select()
